$wb = $excel.ActiveWorkbook

# Sheets "展览" and "全部类型" both carry the same rows of event data; the
# "想去人数" (F column) counts got refreshed for four events.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F3").Value = 1712
    $ws.Range("F4").Value = 31
    $ws.Range("F6").Value = 475
    $ws.Range("F9").Value = 635
}
